$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036566368707438
$ws.Range("D2").Value = 1.044708188474221
$ws.Range("E2").Value = 1.035470784259076
$ws.Range("F2").Value = 1.051816545235598
$ws.Range("I2").Value = 1.039149228070396
$ws.Range("J2").Value = 1.041673986803223
$ws.Range("K2").Value = 1.047478479059697
$ws.Range("L2").Value = 1.038267295479748
$ws.Range("M2").Value = 1.054566998531514

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037442711288285
$ws.Range("D3").Value = 1.045414660823857
$ws.Range("E3").Value = 1.03621411744198
$ws.Range("F3").Value = 1.052700979514288
$ws.Range("I3").Value = 1.039357220548455
$ws.Range("J3").Value = 1.042194715197283
$ws.Range("K3").Value = 1.047996517510298
$ws.Range("L3").Value = 1.038820205951724
$ws.Range("M3").Value = 1.055263970268134

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038010380233053
$ws.Range("D4").Value = 1.045872337216553
$ws.Range("E4").Value = 1.036696005294813
$ws.Range("F4").Value = 1.053274276883259
$ws.Range("I4").Value = 1.039490919280576
$ws.Range("J4").Value = 1.042531628192384
$ws.Range("K4").Value = 1.048331589507334
$ws.Range("L4").Value = 1.039178211139985
$ws.Range("M4").Value = 1.055715327550684

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038249174445168
$ws.Range("D5").Value = 1.046064872349071
$ws.Range("E5").Value = 1.036898805099212
$ws.Range("F5").Value = 1.053515530940069
$ws.Range("I5").Value = 1.039546913373726
$ws.Range("J5").Value = 1.042673257296445
$ws.Range("K5").Value = 1.048472420551438
$ws.Range("L5").Value = 1.039328771524984
$ws.Range("M5").Value = 1.055905165208877

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038289277574571
$ws.Range("D6").Value = 1.046097207301116
$ws.Range("E6").Value = 1.036932868587551
$ws.Range("F6").Value = 1.053556052570748
$ws.Range("I6").Value = 1.039556302527078
$ws.Range("J6").Value = 1.042697036876726
$ws.Range("K6").Value = 1.04849606471863
$ws.Range("L6").Value = 1.039354054464098
$ws.Range("M6").Value = 1.055937044848084

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038013570440392
$ws.Range("D7").Value = 1.04587490937965
$ws.Range("E7").Value = 1.036698714276379
$ws.Range("F7").Value = 1.053277499591578
$ws.Range("I7").Value = 1.039491668313591
$ws.Range("J7").Value = 1.042533520684733
$ws.Range("K7").Value = 1.048333471430759
$ws.Range("L7").Value = 1.039180222720708
$ws.Range("M7").Value = 1.055717863830587

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036862404221456
$ws.Range("D8").Value = 1.044946830979083
$ws.Range("E8").Value = 1.035721809495695
$ws.Range("F8").Value = 1.052115233902202
$ws.Range("I8").Value = 1.039219703215158
$ws.Range("J8").Value = 1.041849975535468
$ws.Range("K8").Value = 1.047653579469023
$ws.Range("L8").Value = 1.038454104699711
$ws.Range("M8").Value = 1.054802465879298

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03483868964502
$ws.Range("D9").Value = 1.043315659625393
$ws.Range("E9").Value = 1.034007358122027
$ws.Range("F9").Value = 1.050074981385657
$ws.Range("I9").Value = 1.038733707565934
$ws.Range("J9").Value = 1.040645281070026
$ws.Range("K9").Value = 1.04645455601858
$ws.Range("L9").Value = 1.037176449884883
$ws.Range("M9").Value = 1.053192318171863

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033492851511073
$ws.Range("D10").Value = 1.042231149721477
$ws.Range("E10").Value = 1.032869183910538
$ws.Range("F10").Value = 1.048720168446303
$ws.Range("I10").Value = 1.038405206446973
$ws.Range("J10").Value = 1.039842087741853
$ws.Range("K10").Value = 1.045654631868969
$ws.Range("L10").Value = 1.036326005943199
$ws.Range("M10").Value = 1.052120926600119

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032910890642273
$ws.Range("D11").Value = 1.041762262701224
$ws.Range("E11").Value = 1.032377500064988
$ws.Range("F11").Value = 1.048134811959934
$ws.Range("I11").Value = 1.03826190193079
$ws.Range("J11").Value = 1.039494296016713
$ws.Range("K11").Value = 1.045308135218164
$ws.Range("L11").Value = 1.03595808498586
$ws.Range("M11").Value = 1.051657505102746

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032694845474547
$ws.Range("D12").Value = 1.041588205984267
$ws.Range("E12").Value = 1.032195041750406
$ws.Range("F12").Value = 1.047917579565778
$ws.Range("I12").Value = 1.038208513453191
$ws.Range("J12").Value = 1.039365111321367
$ws.Range("K12").Value = 1.045179413532065
$ws.Range("L12").Value = 1.035821473169178
$ws.Range("M12").Value = 1.051485446083825

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032741182397373
$ws.Range("D13").Value = 1.041625536790355
$ws.Range("E13").Value = 1.032234171733469
$ws.Range("F13").Value = 1.047964167795712
$ws.Range("I13").Value = 1.038219972638293
$ws.Range("J13").Value = 1.039392821829034
$ws.Range("K13").Value = 1.045207025537939
$ws.Range("L13").Value = 1.035850774565655
$ws.Range("M13").Value = 1.051522349847598

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03289302980443
$ws.Range("D14").Value = 1.041747872889449
$ws.Range("E14").Value = 1.032362414424767
$ws.Range("F14").Value = 1.04811685147277
$ws.Range("I14").Value = 1.038257492058226
$ws.Range("J14").Value = 1.039483617549591
$ws.Range("K14").Value = 1.045297495391146
$ws.Range("L14").Value = 1.035946791574648
$ws.Range("M14").Value = 1.051643281084061

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032986604081536
$ws.Range("D15").Value = 1.041823262665262
$ws.Range("E15").Value = 1.032441452214403
$ws.Range("F15").Value = 1.048210950827427
$ws.Range("I15").Value = 1.038280587988381
$ws.Range("J15").Value = 1.039539559909805
$ws.Range("K15").Value = 1.045353234585327
$ws.Range("L15").Value = 1.036005957541313
$ws.Range("M15").Value = 1.051717800973936

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03353149122864
$ws.Range("D16").Value = 1.042262283373601
$ws.Range("E16").Value = 1.032901839789164
$ws.Range("F16").Value = 1.048759043885599
$ws.Range("I16").Value = 1.038414694790168
$ws.Range("J16").Value = 1.039865169556936
$ws.Range("K16").Value = 1.045677625200952
$ws.Range("L16").Value = 1.036350430643252
$ws.Range("M16").Value = 1.052151692983397

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033873498641351
$ws.Range("D17").Value = 1.042537861529027
$ws.Range("E17").Value = 1.033190938681769
$ws.Range("F17").Value = 1.04910319386634
$ws.Range("I17").Value = 1.038498532750569
$ws.Range("J17").Value = 1.040069415734831
$ws.Range("K17").Value = 1.045881074496794
$ws.Range("L17").Value = 1.036566597951122
$ws.Range("M17").Value = 1.052423996413138

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034073062477745
$ws.Range("D18").Value = 1.042698670278223
$ws.Range("E18").Value = 1.033359676316823
$ws.Range("F18").Value = 1.049304054731695
$ws.Range("I18").Value = 1.038547331563944
$ws.Range("J18").Value = 1.040188548636308
$ws.Range("K18").Value = 1.045999730959444
$ws.Range("L18").Value = 1.036692716136778
$ws.Range("M18").Value = 1.052582874356244

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034141121554722
$ws.Range("D19").Value = 1.042753513502814
$ws.Range("E19").Value = 1.03341723028917
$ws.Range("F19").Value = 1.049372564074823
$ws.Range("I19").Value = 1.038563953300857
$ws.Range("J19").Value = 1.040229169766334
$ws.Range("K19").Value = 1.046040187661344
$ws.Range("L19").Value = 1.036735724474254
$ws.Range("M19").Value = 1.05263705571757

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033836796539025
$ws.Range("D20").Value = 1.042508287485564
$ws.Range("E20").Value = 1.033159909598244
$ws.Range("F20").Value = 1.049066256991141
$ws.Range("I20").Value = 1.03848954832304
$ws.Range("J20").Value = 1.040047502089641
$ws.Range("K20").Value = 1.045859247552088
$ws.Range("L20").Value = 1.036543401982149
$ws.Range("M20").Value = 1.052394775860515

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032848311175735
$ws.Range("D21").Value = 1.041711844940437
$ws.Range("E21").Value = 1.03232464530506
$ws.Range("F21").Value = 1.048071884536188
$ws.Range("I21").Value = 1.038246447902552
$ws.Range("J21").Value = 1.039456880445731
$ws.Range("K21").Value = 1.045270854744543
$ws.Range("L21").Value = 1.03591851555671
$ws.Range("M21").Value = 1.051607667728943

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032227512038296
$ws.Range("D22").Value = 1.041211720218887
$ws.Range("E22").Value = 1.031800494594425
$ws.Range("F22").Value = 1.047447812598148
$ws.Range("I22").Value = 1.038092682581326
$ws.Range("J22").Value = 1.039085537453165
$ws.Range("K22").Value = 1.044900808773725
$ws.Range("L22").Value = 1.035525916950625
$ws.Range("M22").Value = 1.051113224150612

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032556542495643
$ws.Range("D23").Value = 1.041476785485234
$ws.Range("E23").Value = 1.032078260199647
$ws.Range("F23").Value = 1.047778537278768
$ws.Range("I23").Value = 1.038174283324342
$ws.Range("J23").Value = 1.039282392594908
$ws.Range("K23").Value = 1.045096986207888
$ws.Range("L23").Value = 1.035734012760481
$ws.Range("M23").Value = 1.051375295590095

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033853380402412
$ws.Range("D24").Value = 1.042521650510825
$ws.Range("E24").Value = 1.033173929962263
$ws.Range("F24").Value = 1.04908294679218
$ws.Range("I24").Value = 1.03849360831574
$ws.Range("J24").Value = 1.040057403924135
$ws.Range("K24").Value = 1.045869110246313
$ws.Range("L24").Value = 1.036553883145044
$ws.Range("M24").Value = 1.052407979222456

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035361292181416
$ws.Range("D25").Value = 1.043736846016449
$ws.Range("E25").Value = 1.034449747398686
$ws.Range("F25").Value = 1.050601499381006
$ws.Range("I25").Value = 1.038860145888054
$ws.Range("J25").Value = 1.040956739906127
$ws.Range("K25").Value = 1.046764639191818
$ws.Range("L25").Value = 1.037506526197506
$ws.Range("M25").Value = 1.053608226949513
